# Editing existing test cases: rework the SQL Queries section at the
# top of the sheet (rows 2-6), adjust column E width, and adjust row
# heights for a couple of test-case rows further down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. SQL Queries section
# ---------------------------------------------------------------------
# Row 3: new "3rd highest salary" SELECT query (was the old MAX() query)
$ws.Range("B3").Value = "select Employee.EmpName, Salary.salary from Employee Inner join Salary on Employee.EmpID=Salary.EmpID order by salary desc limit 2,1;"

# Row 4: UPDATE query - now references Salary.salary (lower-case) instead of Salary.Salary
$ws.Range("B4").Value = "update Salary inner join Employee on Employee.EmpID=Salary.EmpID set Salary.salary=5000 where datediff(sysdate(),Employee.Date_of_Birth)/365 >30;"

# Row 3 & 4 (A + B) get bumped from the default 11pt to 12pt Times New Roman
# (font/name/bold already match, only the size actually changes)
$ws.Range("A3:B4").Font.Size = 12

# New row 6: "Note:" label (re-uses the existing bold 14pt header style from A2/A7)
# plus explanatory text about the queries (12pt, like the new row 3/4 style).
$ws.Range("A6").Value = "Note:"
$ws.Range("A2").Copy()
$ws.Range("A6").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("B6").Value = "The Query Number 1 will provide you 3rd highest salary along Employee Name from Mysql Database and in the similar way if you want to fetch 2nd highest salary change limit 1,1 and for 4th highest salary change limit 3,1 and so on…. "
$ws.Range("B6").Font.Size = 12

# ---------------------------------------------------------------------
# 2. Layout tweaks
# ---------------------------------------------------------------------
# Column E is narrower now
$ws.Columns.Item(5).ColumnWidth = 51

# Row heights grow for TC_07 (row 17) and TC_10 (row 20) because column E
# got narrower, causing their wrapped text to take up more vertical space.
$ws.Rows.Item(17).RowHeight = 93
$ws.Rows.Item(20).RowHeight = 83.4

# ---------------------------------------------------------------------
# 3. Selection the author ended on
# ---------------------------------------------------------------------
$ws.Range("E14").Select()
